$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordered list of Category values for rows 2..21 (A column holds serial 1..20)
$categories = @(
    "Fan",
    "Light",
    "Cloth",
    "Sport",
    "Gym",
    "Mobile",
    "Door",
    "House",
    "Edu",
    "Fish",
    "Rice",
    "Egg",
    "Fruit",
    "Oil",
    "Deo",
    "Islamic",
    "Sad",
    "Happy",
    "Good",
    "Ram"
)

for ($i = 0; $i -lt $categories.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $categories[$i]
}

$ws.Range("B21").Select()
